# MOS-23045: Update Master Data as per 22 April Changes
# Adds 10 new "Postal Code" location rows (rows 110-119) under the
# "BNMR" (Ben Mansour) parent location, in eng/fra/ara, mirroring the
# existing Postal Code rows already present in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(10110, 10110, 5, "Postal Code",     "BNMR", "eng"),
    @(10111, 10111, 5, "Postal Code",     "BNMR", "eng"),
    @(10113, 10113, 5, "Postal Code",     "BNMR", "eng"),
    @(10114, 10114, 5, "Postal Code",     "BNMR", "eng"),
    @(10111, 10111, 5, "code postal",     "BNMR", "fra"),
    @(10110, 10110, 5, "code postal",     "BNMR", "fra"),
    @(10113, 10113, 5, "code postal",     "BNMR", "fra"),
    @(10114, 10114, 5, "code postal",     "BNMR", "fra"),
    @(10111, 10111, 5, "الرمز البريدي",    "BNMR", "ara"),
    @(10110, 10110, 5, "الرمز البريدي",    "BNMR", "ara")
)

$startRow = 110
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $true
    $ws.Cells.Item($r, 8).Value = "superadmin"
    $ws.Cells.Item($r, 9).Value = "now()"
}
